$wb = $excel.ActiveWorkbook

# The same update needs to be applied to both the "展览" sheet and the
# "全部类型" sheet, which mirror each other's data.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 1053
    $ws.Range("F5").Value = 73

    $ws.Range("F6").Value = 3020
    $ws.Range("G6").Value = "不可售"

    $ws.Range("F8").Value = 2094
    $ws.Range("F9").Value = 186
    $ws.Range("F11").Value = 943
    $ws.Range("F13").Value = 33
    $ws.Range("F14").Value = 232
    $ws.Range("F15").Value = 89
    $ws.Range("F17").Value = 34
}
